$wb = $excel.ActiveWorkbook

# ---- Step 1: Add sheet "2301260924" as an exact copy of the last sheet (snapshot) ----
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Copy($null, $lastSheet)
$sheet5 = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet5.Name = "2301260924"

# ---- Step 2: Add sheet "2301260928" as a copy of sheet5, then apply edits ----
$sheet5.Copy($null, $sheet5)
$sheet6 = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet6.Name = "2301260928"

# Clear old data rows (2-51), keep header row 1 intact
$sheet6.Range("A2:G51").ClearContents()

$sheet6.Cells.Item(2,1).Value = "fgb-acesh"
$sheet6.Cells.Item(2,2).Value = "Academic Center for Education, Sport and Health"
$sheet6.Cells.Item(2,3).Value = "vu-alab"
$sheet6.Cells.Item(2,4).Value = "A-LAB"
$sheet6.Cells.Item(2,5).Value = "acta"
$sheet6.Cells.Item(2,6).Value = "Academic Centre for Dentistry Amsterdam"
$sheet6.Cells.Item(2,7).Value = "acta-fda"
$sheet6.Cells.Item(3,1).Value = "sbe-acc"
$sheet6.Cells.Item(3,2).Value = "Accounting"
$sheet6.Cells.Item(3,3).Value = "vu-aimms"
$sheet6.Cells.Item(3,4).Value = "AIMMS"
$sheet6.Cells.Item(3,5).Value = "fgb"
$sheet6.Cells.Item(3,6).Value = "Faculty of Behavioural and Movement Sciences"
$sheet6.Cells.Item(3,7).Value = "vu-cic"
$sheet6.Cells.Item(4,1).Value = "beta-aile"
$sheet6.Cells.Item(4,2).Value = "Amsterdam Institute for Life and Environment"
$sheet6.Cells.Item(4,3).Value = "vu-abri"
$sheet6.Cells.Item(4,4).Value = "Amsterdam Business Research Institute"
$sheet6.Cells.Item(4,5).Value = "fgw"
$sheet6.Cells.Item(4,6).Value = "Faculty of Humanities"
$sheet6.Cells.Item(4,7).Value = "beta-ecsc"
$sheet6.Cells.Item(5,1).Value = "fgw-acha"
$sheet6.Cells.Item(5,2).Value = "Art and Culture, History, Antiquity"
$sheet6.Cells.Item(5,3).Value = "vu-ams"
$sheet6.Cells.Item(5,4).Value = "Amsterdam Movement Sciences"
$sheet6.Cells.Item(5,5).Value = "rch"
$sheet6.Cells.Item(5,6).Value = "Faculty of Law"
$sheet6.Cells.Item(5,7).Value = "vu-acwfs"
$sheet6.Cells.Item(6,1).Value = "beta-ai"
$sheet6.Cells.Item(6,2).Value = "Athena Institute"
$sheet6.Cells.Item(6,3).Value = "vu-an"
$sheet6.Cells.Item(6,4).Value = "Amsterdam Neuroscience"
$sheet6.Cells.Item(6,5).Value = "frt"
$sheet6.Cells.Item(6,6).Value = "Faculty of Religion and Theology"
$sheet6.Cells.Item(6,7).Value = "acta-acdia"
$sheet6.Cells.Item(7,1).Value = "frt-bp"
$sheet6.Cells.Item(7,2).Value = "Beliefs and Practices"
$sheet6.Cells.Item(7,3).Value = "vu-aph"
$sheet6.Cells.Item(7,4).Value = "Amsterdam Public Health"
$sheet6.Cells.Item(7,5).Value = "beta"
$sheet6.Cells.Item(7,6).Value = "Faculty of Science"
$sheet6.Cells.Item(7,7).Value = "vu-whocc"
$sheet6.Cells.Item(8,1).Value = "fgb-bp"
$sheet6.Cells.Item(8,2).Value = "Biological Psychology"
$sheet6.Cells.Item(8,3).Value = "vu-ard"
$sheet6.Cells.Item(8,4).Value = "Amsterdam Reproduction & Development"
$sheet6.Cells.Item(8,5).Value = "fsw"
$sheet6.Cells.Item(8,6).Value = "Faculty of Social Sciences"
$sheet6.Cells.Item(8,7).Value = "vu-kcdi"
$sheet6.Cells.Item(9,1).Value = "beta-cncr"
$sheet6.Cells.Item(9,2).Value = "Center for Neurogenomics and Cognitive Research"
$sheet6.Cells.Item(9,3).Value = "vu-asi"
$sheet6.Cells.Item(9,4).Value = "Amsterdam Sustainability Institute"
$sheet6.Cells.Item(9,5).Value = "sbe"
$sheet6.Cells.Item(9,6).Value = "School of Business and Economics"
$sheet6.Cells.Item(9,7).Value = "beta-mcb"
$sheet6.Cells.Item(10,1).Value = "beta-cps"
$sheet6.Cells.Item(10,2).Value = "Chemistry and Pharmaceutical Sciences"
$sheet6.Cells.Item(10,3).Value = "vu-clue"
$sheet6.Cells.Item(10,4).Value = "CLUE+"
$sheet6.Cells.Item(10,5).Value = "gnk"
$sheet6.Cells.Item(10,6).Value = "VUmc - School of Medical Sciences"
$sheet6.Cells.Item(10,7).Value = "beta-eh"
$sheet6.Cells.Item(11,1).Value = "fgb-cndp"
$sheet6.Cells.Item(11,2).Value = "Clinical, Neuro- & Developmental Psychology"
$sheet6.Cells.Item(11,3).Value = "vu-ibba"
$sheet6.Cells.Item(11,4).Value = "IBBA"
$sheet6.Cells.Item(12,1).Value = "fsw-cosc"
$sheet6.Cells.Item(12,2).Value = "Communication Science"
$sheet6.Cells.Item(12,3).Value = "vu-isr"
$sheet6.Cells.Item(12,4).Value = "Institute for Societal Resilience"
$sheet6.Cells.Item(13,1).Value = "beta-cs"
$sheet6.Cells.Item(13,2).Value = "Computer Science"
$sheet6.Cells.Item(13,3).Value = "vu-ki"
$sheet6.Cells.Item(13,4).Value = "Kooijmans Institute"
$sheet6.Cells.Item(14,1).Value = "rch-cal"
$sheet6.Cells.Item(14,2).Value = "Constitutional and Administrative Law"
$sheet6.Cells.Item(14,3).Value = "vu-learn"
$sheet6.Cells.Item(14,4).Value = "LEARN!"
$sheet6.Cells.Item(15,1).Value = "rch-clc"
$sheet6.Cells.Item(15,2).Value = "Criminal Law and Criminology"
$sheet6.Cells.Item(15,3).Value = "vu-laser"
$sheet6.Cells.Item(15,4).Value = "LaserLaB"
$sheet6.Cells.Item(16,1).Value = "rch-dpl"
$sheet6.Cells.Item(16,2).Value = "Dutch Private Law"
$sheet6.Cells.Item(16,3).Value = "vu-ni"
$sheet6.Cells.Item(16,4).Value = "Network Institute"
$sheet6.Cells.Item(17,1).Value = "beta-es"
$sheet6.Cells.Item(17,2).Value = "Earth Sciences"
$sheet6.Cells.Item(17,3).Value = "vu-tain"
$sheet6.Cells.Item(17,4).Value = "Talma Institute"
$sheet6.Cells.Item(18,1).Value = "sbe-eds"
$sheet6.Cells.Item(18,2).Value = "Econometrics and Data Science"
$sheet6.Cells.Item(18,3).Value = "vu-ti"
$sheet6.Cells.Item(18,4).Value = "Tinbergen Institute"
$sheet6.Cells.Item(19,1).Value = "sbe-econ"
$sheet6.Cells.Item(19,2).Value = "Economics"
$sheet6.Cells.Item(20,1).Value = "fgb-efs"
$sheet6.Cells.Item(20,2).Value = "Educational and Family Studies"
$sheet6.Cells.Item(21,1).Value = "sbe-egs"
$sheet6.Cells.Item(21,2).Value = "Ethics, Governance and Society"
$sheet6.Cells.Item(22,1).Value = "fgb-eap"
$sheet6.Cells.Item(22,2).Value = "Experimental and Applied Psychology"
$sheet6.Cells.Item(23,1).Value = "sbe-fin"
$sheet6.Cells.Item(23,2).Value = "Finance"
$sheet6.Cells.Item(24,1).Value = "beta-hs"
$sheet6.Cells.Item(24,2).Value = "Health Sciences"
$sheet6.Cells.Item(25,1).Value = "beta-hsas"
$sheet6.Cells.Item(25,2).Value = "History and Social Aspects of Science"
$sheet6.Cells.Item(26,1).Value = "fgb-hms"
$sheet6.Cells.Item(26,2).Value = "Human Movement Sciences"
$sheet6.Cells.Item(27,1).Value = "beta-ies"
$sheet6.Cells.Item(27,2).Value = "Institute for Environmental Studies"
$sheet6.Cells.Item(28,1).Value = "beta-kggb"
$sheet6.Cells.Item(28,2).Value = "Kars Group (Geo- and Bioarchaeology)"
$sheet6.Cells.Item(29,1).Value = "sbe-kii"
$sheet6.Cells.Item(29,2).Value = "Knowledge, Information and Innovation"
$sheet6.Cells.Item(30,1).Value = "fgw-llc"
$sheet6.Cells.Item(30,2).Value = "Language, Literature and Communication"
$sheet6.Cells.Item(31,1).Value = "rch-ltlh"
$sheet6.Cells.Item(31,2).Value = "Legal Theory and Legal History"
$sheet6.Cells.Item(32,1).Value = "sbe-mo"
$sheet6.Cells.Item(32,2).Value = "Management and Organisation"
$sheet6.Cells.Item(33,1).Value = "sbe-mrk"
$sheet6.Cells.Item(33,2).Value = "Marketing"
$sheet6.Cells.Item(34,1).Value = "beta-math"
$sheet6.Cells.Item(34,2).Value = "Mathematics"
$sheet6.Cells.Item(35,1).Value = "rch-ntl"
$sheet6.Cells.Item(35,2).Value = "Notary and Tax Law"
$sheet6.Cells.Item(36,1).Value = "acta-owi"
$sheet6.Cells.Item(36,2).Value = "OWI (ACTA)"
$sheet6.Cells.Item(37,1).Value = "sbe-oa"
$sheet6.Cells.Item(37,2).Value = "Operations Analytics"
$sheet6.Cells.Item(38,1).Value = "acta-oii"
$sheet6.Cells.Item(38,2).Value = "Oral Infections and Inflammation (OII)"
$sheet6.Cells.Item(39,1).Value = "acta-orm"
$sheet6.Cells.Item(39,2).Value = "Oral Regenerative Medicine (ORM)"
$sheet6.Cells.Item(40,1).Value = "fsw-os"
$sheet6.Cells.Item(40,2).Value = "Organization Sciences"
$sheet6.Cells.Item(41,1).Value = "fgw-phil"
$sheet6.Cells.Item(41,2).Value = "Philosophy"
$sheet6.Cells.Item(42,1).Value = "beta-pa"
$sheet6.Cells.Item(42,2).Value = "Physics and Astronomy"
$sheet6.Cells.Item(43,1).Value = "fsw-pspa"
$sheet6.Cells.Item(43,2).Value = "Political Science and Public Administration"
$sheet6.Cells.Item(44,1).Value = "fsw-sca"
$sheet6.Cells.Item(44,2).Value = "Social and Cultural Anthropology"
$sheet6.Cells.Item(45,1).Value = "fsw-socio"
$sheet6.Cells.Item(45,2).Value = "Sociology"
$sheet6.Cells.Item(46,1).Value = "sbe-se"
$sheet6.Cells.Item(46,2).Value = "Spatial Economics"
$sheet6.Cells.Item(47,1).Value = "frt-tt"
$sheet6.Cells.Item(47,2).Value = "Texts and Traditions"
$sheet6.Cells.Item(48,1).Value = "rch-tls"
$sheet6.Cells.Item(48,2).Value = "Transnational Legal Studies"
$sheet6.Cells.Item(49,1).Value = "sbe-vsee"
$sheet6.Cells.Item(49,2).Value = "VU SBE Executive Education"

# Restore the originally active tab (first sheet), since adding/renaming sheets
# shifts Excel's active-tab selection as a side effect.
$wb.Worksheets.Item(1).Activate()

Write-Host "Done. Sheet count:" $wb.Worksheets.Count
foreach ($s in $wb.Worksheets) {
    Write-Host $s.Name
}
